$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition")

$ws.Range("B4").Value = "Farquat Homeowners Assoc"
$ws.Range("B6").Value = "County of Windsor Assessors Office"
$ws.Range("B8").Value = "Central Gas Electric"
$ws.Range("B9").Value = "Central Gas Electric"
$ws.Range("B14").Value = "Bed Bath Beyond,Target,Container Store"
$ws.Range("B16").Value = "Starbucks,Uptown Espresso,Tim Hortons"
$ws.Range("B32").Value = "Speedy Speeds"
$ws.Range("B33").Value = "Megacorp Inc"
$ws.Range("B28").Value = "Mikes Wrenchems"

$ws.Range("B17").Select()
